# Apply "add almost all lineal" regeneration of the random experiment
# values. All the changed cells hold numbers that were actually written
# as TEXT (shared strings) in the source workbook, so every numeric-
# looking replacement is entered with a leading apostrophe to force
# Excel to keep storing it as text instead of re-parsing it as a
# number (which would also introduce float round-trip noise).

$wb = $excel.ActiveWorkbook

# NOTE: the workbook has two sheets whose names differ only by case
# ("Vector_bf" vs "Vector_BF"); Worksheets.Item(name) resolves
# case-insensitively and would hit the wrong one, so every sheet below
# is addressed by its 1-based tab index instead:
#   1 Funciones_Objetivo        5 Vector_bf
#   2 Restricciones_del_lider   6 Vector_BF
#   3 Restricciones_del_follower 7 Vector_Alpha
#   4 Punto_modificado

# --- Restricciones_del_follower (index 3) -------------------------------
$ws = $wb.Worksheets.Item(3)

$ws.Range("A4").Value = "-16 - 2x + y_1 + 4y_2"

$ws.Range("B2").Value = "'-4.382729079133727"
$ws.Range("B3").Value = "'0.38272907913372656"
$ws.Range("B4").Value = "'-13.81788203899294"
$ws.Range("B5").Value = "'3.3894519012124054"
$ws.Range("B6").Value = "'-8.517138863787123"

$ws.Range("D2").Value = "'0.0866877650392671"
$ws.Range("D3").Value = "'0.9648587319705634"
$ws.Range("D4").Value = "'0.9761226555169311"
$ws.Range("D5").Value = "'0.8143958706897286"
$ws.Range("D6").Value = "'0.23927405565041526"

$ws.Range("E2").Value = "'0.2228861160164336"
$ws.Range("E3").Value = "'0.1484093051209633"
$ws.Range("E4").Value = "'0.10580998854943345"
$ws.Range("E5").Value = "'0.9094608608554801"
$ws.Range("E6").Value = "'0.8112793847527293"

$ws.Range("F2").Value = "'0.14311793845404397"
$ws.Range("F3").Value = "'0.599513580292495"
$ws.Range("F4").Value = "'0.02723005887389917"
$ws.Range("F5").Value = "'0.5401584619592554"
$ws.Range("F6").Value = "'0.8119833452409898"

# --- Punto_modificado (index 4) ------------------------------------------
$ws = $wb.Worksheets.Item(4)

$ws.Range("A2").Value = "'5.875840352759835"
$ws.Range("B2").Value = "'4.382729079133727"
$ws.Range("C2").Value = "'2.387767396848251"

# --- Vector_bf (index 5) ----------------------------------------------------
$ws = $wb.Worksheets.Item(5)

$ws.Range("A2").Value = "'-1.1901413818371251"
$ws.Range("A3").Value = "'-3.9044906220677245"

# --- Vector_BF (index 6) -----------------------------------------------------
$ws = $wb.Worksheets.Item(6)

$ws.Range("A2").Value = "'-4.441508140239515"
$ws.Range("A3").Value = "'3.6817647309960155"
$ws.Range("A4").Value = "'-2.423239954197734"
